$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared-string values introduced by the bug-fix pass
$ws.Range("H9").Value = "Resolved"
$ws.Range("C13").Value = "Snehal/Mukesh"
$ws.Range("H13").Value = "Resolved e,f"
$ws.Range("H14").Value = "Resolved c,d"
$ws.Range("H15").Value = "Resolved"
$ws.Range("H16").Value = "Resolved"

# Update the view state: scrolled position and active cell selection
$ws.Activate()
$app = $ws.Application
$app.ActiveWindow.ScrollRow = 6
$app.ActiveWindow.ScrollColumn = 1
$ws.Range("I9").Select()
